$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.045595554491117
$ws.Cells.Item(2, 4).Value = 1.055080504289956
$ws.Cells.Item(2, 5).Value = 1.043296372145855
$ws.Cells.Item(2, 6).Value = 1.062350228400208
$ws.Cells.Item(2, 9).Value = 1.03810074696228
$ws.Cells.Item(2, 10).Value = 1.050654447703432
$ws.Cells.Item(2, 11).Value = 1.057821945887476
$ws.Cells.Item(2, 12).Value = 1.046070637592752
$ws.Cells.Item(2, 13).Value = 1.065071813958662
$ws.Cells.Item(2, 14).Value = 1.020652616663084

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.046849795356579
$ws.Cells.Item(3, 4).Value = 1.056276311228005
$ws.Cells.Item(3, 5).Value = 1.044372528758724
$ws.Cells.Item(3, 6).Value = 1.063610679061898
$ws.Cells.Item(3, 9).Value = 1.038332419047667
$ws.Cells.Item(3, 10).Value = 1.051554910589394
$ws.Cells.Item(3, 11).Value = 1.05883014729701
$ws.Cells.Item(3, 12).Value = 1.046957106993125
$ws.Cells.Item(3, 13).Value = 1.066145942945814
$ws.Cells.Item(3, 14).Value = 1.020959516143102

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.047660940863848
$ws.Cells.Item(4, 4).Value = 1.057049939287781
$ws.Cells.Item(4, 5).Value = 1.045068800036137
$ws.Cells.Item(4, 6).Value = 1.06442619151544
$ws.Cells.Item(4, 9).Value = 1.038480674095221
$ws.Cells.Item(4, 10).Value = 1.052136639724226
$ws.Cells.Item(4, 11).Value = 1.059481811018407
$ws.Cells.Item(4, 12).Value = 1.047530044948852
$ws.Cells.Item(4, 13).Value = 1.066840332797695
$ws.Cells.Item(4, 14).Value = 1.021157592516638

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.04800184560818
$ws.Cells.Item(5, 4).Value = 1.057375141189229
$ws.Cells.Item(5, 5).Value = 1.045361496267711
$ws.Cells.Item(5, 6).Value = 1.064769015082173
$ws.Cells.Item(5, 9).Value = 1.038542605363985
$ws.Cells.Item(5, 10).Value = 1.052380977461959
$ws.Cells.Item(5, 11).Value = 1.05975560220407
$ws.Cells.Item(5, 12).Value = 1.047770749804343
$ws.Cells.Item(5, 13).Value = 1.067132102092301
$ws.Cells.Item(5, 14).Value = 1.021240742610326

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.048059079179011
$ws.Cells.Item(6, 4).Value = 1.057429742238235
$ws.Cells.Item(6, 5).Value = 1.045410640345805
$ws.Cells.Item(6, 6).Value = 1.064826575680627
$ws.Cells.Item(6, 9).Value = 1.038552980736324
$ws.Cells.Item(6, 10).Value = 1.05242198991862
$ws.Cells.Item(6, 11).Value = 1.059801563133359
$ws.Cells.Item(6, 12).Value = 1.047811155969569
$ws.Cells.Item(6, 13).Value = 1.067181082569035
$ws.Cells.Item(6, 14).Value = 1.021254696781332

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.047665496442611
$ws.Cells.Item(7, 4).Value = 1.057054284772214
$ws.Cells.Item(7, 5).Value = 1.045072711119345
$ws.Cells.Item(7, 6).Value = 1.064430772408314
$ws.Cells.Item(7, 9).Value = 1.038481503176044
$ws.Cells.Item(7, 10).Value = 1.052139905444192
$ws.Cells.Item(7, 11).Value = 1.059485470088554
$ws.Cells.Item(7, 12).Value = 1.047533261879
$ws.Cells.Item(7, 13).Value = 1.066844232029043
$ws.Cells.Item(7, 14).Value = 1.021158704048014

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.046019522146419
$ws.Cells.Item(8, 4).Value = 1.055484662791696
$ws.Cells.Item(8, 5).Value = 1.043660080544179
$ws.Cells.Item(8, 6).Value = 1.062776222138476
$ws.Cells.Item(8, 9).Value = 1.038179383990077
$ws.Cells.Item(8, 10).Value = 1.050958956989145
$ws.Cells.Item(8, 11).Value = 1.058162820331572
$ws.Cells.Item(8, 12).Value = 1.046370363048312
$ws.Cells.Item(8, 13).Value = 1.06543495521349
$ws.Cells.Item(8, 14).Value = 1.020756440246447

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.043115676312219
$ws.Cells.Item(9, 4).Value = 1.052717636453305
$ws.Cells.Item(9, 5).Value = 1.041170210192791
$ws.Cells.Item(9, 6).Value = 1.059859959860765
$ws.Cells.Item(9, 9).Value = 1.037634338346701
$ws.Cells.Item(9, 10).Value = 1.048870776305333
$ws.Cells.Item(9, 11).Value = 1.05582662651401
$ws.Cells.Item(9, 12).Value = 1.044316013812232
$ws.Cells.Item(9, 13).Value = 1.062946618683773
$ws.Cells.Item(9, 14).Value = 1.020043689022816

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.041177305592098
$ws.Cells.Item(10, 4).Value = 1.050872053025046
$ws.Cells.Item(10, 5).Value = 1.039509772223312
$ws.Cells.Item(10, 6).Value = 1.057915157215785
$ws.Cells.Item(10, 9).Value = 1.037262426784607
$ws.Cells.Item(10, 10).Value = 1.0474737135106
$ws.Cells.Item(10, 11).Value = 1.054265343676362
$ws.Cells.Item(10, 12).Value = 1.042942877486504
$ws.Cells.Item(10, 13).Value = 1.061284240417465
$ws.Cells.Item(10, 14).Value = 1.019565863734345

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.040337344879572
$ws.Cells.Item(11, 4).Value = 1.050072653175449
$ws.Cells.Item(11, 5).Value = 1.038790637344924
$ws.Cells.Item(11, 6).Value = 1.057072856488722
$ws.Cells.Item(11, 9).Value = 1.037099351091667
$ws.Cells.Item(11, 10).Value = 1.046867574141601
$ws.Cells.Item(11, 11).Value = 1.053588359244361
$ws.Cells.Item(11, 12).Value = 1.042347426279788
$ws.Cells.Item(11, 13).Value = 1.06056355803543
$ws.Cells.Item(11, 14).Value = 1.01935832310099

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.040025247532708
$ws.Cells.Item(12, 4).Value = 1.049775680350255
$ws.Cells.Item(12, 5).Value = 1.038523493416393
$ws.Cells.Item(12, 6).Value = 1.056759957442644
$ws.Cells.Item(12, 9).Value = 1.037038471082095
$ws.Cells.Item(12, 10).Value = 1.046642244291046
$ws.Cells.Item(12, 11).Value = 1.053336753756346
$ws.Cells.Item(12, 12).Value = 1.042126115943918
$ws.Cells.Item(12, 13).Value = 1.06029573262409
$ws.Cells.Item(12, 14).Value = 1.019281136635381

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.040092197986517
$ws.Cells.Item(13, 4).Value = 1.049839383886475
$ws.Cells.Item(13, 5).Value = 1.038580797869761
$ws.Cells.Item(13, 6).Value = 1.056827076799758
$ws.Cells.Item(13, 9).Value = 1.037051543927701
$ws.Cells.Item(13, 10).Value = 1.046690586610618
$ws.Cells.Item(13, 11).Value = 1.053390730510591
$ws.Cells.Item(13, 12).Value = 1.042173593822765
$ws.Cells.Item(13, 13).Value = 1.060353188091361
$ws.Cells.Item(13, 14).Value = 1.01929769777701

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.040311548845267
$ws.Cells.Item(14, 4).Value = 1.050048106142528
$ws.Cells.Item(14, 5).Value = 1.038768555670655
$ws.Cells.Item(14, 6).Value = 1.057046992805408
$ws.Cells.Item(14, 9).Value = 1.037094324985893
$ws.Cells.Item(14, 10).Value = 1.046848952045194
$ws.Cells.Item(14, 11).Value = 1.05356756439343
$ws.Cells.Item(14, 12).Value = 1.042329135427685
$ws.Cells.Item(14, 13).Value = 1.060541422207956
$ws.Cells.Item(14, 14).Value = 1.019351944818513

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.040446684990216
$ws.Cells.Item(15, 4).Value = 1.050176701419056
$ws.Cells.Item(15, 5).Value = 1.038884236046218
$ws.Cells.Item(15, 6).Value = 1.057182486117235
$ws.Cells.Item(15, 9).Value = 1.037120643183394
$ws.Cells.Item(15, 10).Value = 1.046946501942983
$ws.Cells.Item(15, 11).Value = 1.053676498523548
$ws.Cells.Item(15, 12).Value = 1.04242495203014
$ws.Cells.Item(15, 13).Value = 1.060657381910399
$ws.Cells.Item(15, 14).Value = 1.019385355387904

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.041233037385202
$ws.Cells.Item(16, 4).Value = 1.050925101073775
$ws.Cells.Item(16, 5).Value = 1.039557495412703
$ws.Cells.Item(16, 6).Value = 1.057971053737793
$ws.Cells.Item(16, 9).Value = 1.037273206642061
$ws.Cells.Item(16, 10).Value = 1.047513915479638
$ws.Cells.Item(16, 11).Value = 1.054310252930595
$ws.Cells.Item(16, 12).Value = 1.042982377041807
$ws.Cells.Item(16, 13).Value = 1.061332051372473
$ws.Cells.Item(16, 14).Value = 1.019579624004819

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.041726123506864
$ws.Cells.Item(17, 4).Value = 1.051394483655177
$ws.Cells.Item(17, 5).Value = 1.039979770486008
$ws.Cells.Item(17, 6).Value = 1.058465649004185
$ws.Cells.Item(17, 9).Value = 1.03736836026206
$ws.Cells.Item(17, 10).Value = 1.047869515707728
$ws.Cells.Item(17, 11).Value = 1.054707537396634
$ws.Cells.Item(17, 12).Value = 1.043331799646544
$ws.Cells.Item(17, 13).Value = 1.061755021302886
$ws.Cells.Item(17, 14).Value = 1.019701312031772

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.042013671415293
$ws.Cells.Item(18, 4).Value = 1.051668242622528
$ws.Cells.Item(18, 5).Value = 1.040226061645347
$ws.Cells.Item(18, 6).Value = 1.058754120070272
$ws.Cells.Item(18, 9).Value = 1.03742366540195
$ws.Cells.Item(18, 10).Value = 1.048076815477482
$ws.Cells.Item(18, 11).Value = 1.054939176294271
$ws.Cells.Item(18, 12).Value = 1.043535527731662
$ws.Cells.Item(18, 13).Value = 1.06200164953941
$ws.Cells.Item(18, 14).Value = 1.019772228974212

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.042111707595738
$ws.Cells.Item(19, 4).Value = 1.051761583370149
$ws.Cells.Item(19, 5).Value = 1.040310038125666
$ws.Cells.Item(19, 6).Value = 1.058852478281186
$ws.Cells.Item(19, 9).Value = 1.03744248973846
$ws.Cells.Item(19, 10).Value = 1.048147479759175
$ws.Cells.Item(19, 11).Value = 1.055018143906317
$ws.Cells.Item(19, 12).Value = 1.043604979518194
$ws.Cells.Item(19, 13).Value = 1.062085729394962
$ws.Cells.Item(19, 14).Value = 1.019796399384559

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.041673226355783
$ws.Cells.Item(20, 4).Value = 1.051344125843449
$ws.Cells.Item(20, 5).Value = 1.039934465908078
$ws.Cells.Item(20, 6).Value = 1.05841258549433
$ws.Cells.Item(20, 9).Value = 1.037358171494127
$ws.Cells.Item(20, 10).Value = 1.047831375165456
$ws.Cells.Item(20, 11).Value = 1.054664921902527
$ws.Cells.Item(20, 12).Value = 1.043294318632836
$ws.Cells.Item(20, 13).Value = 1.061709649232788
$ws.Cells.Item(20, 14).Value = 1.019688262437033

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.040246958231887
$ws.Cells.Item(21, 4).Value = 1.049986643768735
$ws.Cells.Item(21, 5).Value = 1.038713266390503
$ws.Cells.Item(21, 6).Value = 1.056982233902398
$ws.Cells.Item(21, 9).Value = 1.037081735494688
$ws.Cells.Item(21, 10).Value = 1.046802322432367
$ws.Cells.Item(21, 11).Value = 1.053515495197492
$ws.Cells.Item(21, 12).Value = 1.042283335995488
$ws.Cells.Item(21, 13).Value = 1.060485995612565
$ws.Cells.Item(21, 14).Value = 1.019335973089356

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.039349635056609
$ws.Cells.Item(22, 4).Value = 1.049132907466586
$ws.Cells.Item(22, 5).Value = 1.037945302880993
$ws.Cells.Item(22, 6).Value = 1.056082734052319
$ws.Cells.Item(22, 9).Value = 1.036906155954492
$ws.Cells.Item(22, 10).Value = 1.046154258030346
$ws.Cells.Item(22, 11).Value = 1.052791974611749
$ws.Cells.Item(22, 12).Value = 1.041646919212635
$ws.Cells.Item(22, 13).Value = 1.059715871746943
$ws.Cells.Item(22, 14).Value = 1.019113915395933

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.039825378381005
$ws.Cells.Item(23, 4).Value = 1.049585512219907
$ws.Cells.Item(23, 5).Value = 1.038352429475556
$ws.Cells.Item(23, 6).Value = 1.056559594025403
$ws.Cells.Item(23, 9).Value = 1.036999402283237
$ws.Cells.Item(23, 10).Value = 1.046497910238512
$ws.Cells.Item(23, 11).Value = 1.053175606048075
$ws.Cells.Item(23, 12).Value = 1.04198436959886
$ws.Cells.Item(23, 13).Value = 1.060124202240098
$ws.Cells.Item(23, 14).Value = 1.019231685610553

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.041697128489467
$ws.Cells.Item(24, 4).Value = 1.051366880444946
$ws.Cells.Item(24, 5).Value = 1.039954937142327
$ws.Cells.Item(24, 6).Value = 1.058436562665902
$ws.Cells.Item(24, 9).Value = 1.037362775966514
$ws.Cells.Item(24, 10).Value = 1.047848609592644
$ws.Cells.Item(24, 11).Value = 1.054684178289294
$ws.Cells.Item(24, 12).Value = 1.043311254949775
$ws.Cells.Item(24, 13).Value = 1.061730151174594
$ws.Cells.Item(24, 14).Value = 1.019694159177608

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.043866813676604
$ws.Cells.Item(25, 4).Value = 1.053433128397861
$ws.Cells.Item(25, 5).Value = 1.041813986601428
$ws.Cells.Item(25, 6).Value = 1.060613984301423
$ws.Cells.Item(25, 9).Value = 1.037776750314075
$ws.Cells.Item(25, 10).Value = 1.04941148464994
$ws.Cells.Item(25, 11).Value = 1.056431253748177
$ws.Cells.Item(25, 12).Value = 1.044847734810727
$ws.Cells.Item(25, 13).Value = 1.063590519266368
$ws.Cells.Item(25, 14).Value = 1.020228418669272
